# fix the inference bug, export quality check and design change (variation per task)
#
# Behavioural change: the QualityRate column (D2:D11) on the "WorkMethod"
# sheet is bumped from 0.7 to 0.9 for every work method. The user made this
# edit with the WorkMethod sheet active and the D2:D11 range selected, which
# is why it becomes the new active tab / active selection in the saved file
# (previously "Project" / A2 was active).

$wb = $excel.ActiveWorkbook

# --- Update the QualityRate values on WorkMethod (D2:D11): 0.7 -> 0.9 ---
$wsWorkMethod = $wb.Worksheets.Item("WorkMethod")
$wsWorkMethod.Range("D2:D11").Value = 0.9

# --- Make WorkMethod the active sheet with D2:D11 selected (D2 active) ---
# (leaves the other sheets' own selections/active cells untouched, matching
# the fact that only the WorkMethod sheet's selection genuinely changed)
$wsWorkMethod.Activate()
$wsWorkMethod.Range("D2:D11").Select()
